$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update report header text (volume number + week-covering dates) ---
$ws.Range("A8").Value = "Volume 29   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/5/2022  Through  12/11/2022"

# --- Precinct crime-stat table updates (rows 14-30) ---
$ws.Range("J14").Copy($ws.Range("C14"))
$ws.Range("C14").Value = 1
$ws.Range("J14").Copy($ws.Range("F14"))
$ws.Range("F14").Value = 1
$ws.Range("I14").Value = 4
$ws.Range("K14").Value = -20
$ws.Range("L14").Value = 300
$ws.Range("M14").Value = 300
$ws.Range("N14").Value = -50
$ws.Range("J14").Copy($ws.Range("D15"))
$ws.Range("D15").Value = 1
$ws.Range("H15").Copy($ws.Range("E15"))
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 2
$ws.Range("J15").Value = 17
$ws.Range("K15").Value = -5.882352941176
$ws.Range("L15").Value = 14.285714285714
$ws.Range("N15").Value = 45.454545454545
$ws.Range("J14").Copy($ws.Range("C16"))
$ws.Range("C16").Value = 2
$ws.Range("G14").Copy($ws.Range("D16"))
$ws.Range("H14").Copy($ws.Range("E16"))
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 83.333333333333
$ws.Range("I16").Value = 96
$ws.Range("K16").Value = 39.130434782608
$ws.Range("L16").Value = 54.838709677419
$ws.Range("M16").Value = -2.040816326530
$ws.Range("N16").Value = -87.739463601532
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = 45.454545454545
$ws.Range("I17").Value = 202
$ws.Range("J17").Value = 145
$ws.Range("K17").Value = 39.310344827586
$ws.Range("L17").Value = 26.25
$ws.Range("M17").Value = 56.589147286821
$ws.Range("N17").Value = -29.861111111111
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 165
$ws.Range("J18").Value = 106
$ws.Range("K18").Value = 55.660377358490
$ws.Range("L18").Value = -9.340659340659
$ws.Range("M18").Value = 36.363636363636
$ws.Range("N18").Value = -77.489768076398
$ws.Range("C19").Value = 19
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = 18.75
$ws.Range("F19").Value = 51
$ws.Range("G19").Value = 52
$ws.Range("H19").Value = -1.923076923076
$ws.Range("I19").Value = 596
$ws.Range("J19").Value = 400
$ws.Range("K19").Value = 49
$ws.Range("L19").Value = 80.060422960725
$ws.Range("M19").Value = 6.618962432915
$ws.Range("N19").Value = -60.345974717232
$ws.Range("C20").Value = 1
$ws.Range("J14").Copy($ws.Range("D20"))
$ws.Range("D20").Value = 2
$ws.Range("H15").Copy($ws.Range("E20"))
$ws.Range("E20").Value = -50
$ws.Range("I20").Value = 29
$ws.Range("J20").Value = 31
$ws.Range("K20").Value = -6.451612903225
$ws.Range("L20").Value = 26.086956521739
$ws.Range("M20").Value = 107.142857142857
$ws.Range("N20").Value = -90.9375
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = 22.727272727272
$ws.Range("F21").Value = 95
$ws.Range("G21").Value = 87
$ws.Range("H21").Value = 9.195402298850
$ws.Range("I21").Value = 1108
$ws.Range("J21").Value = 773
$ws.Range("K21").Value = 43.337645536869
$ws.Range("L21").Value = 43.337645536869
$ws.Range("M21").Value = 18.884120171673
$ws.Range("N21").Value = -69.610532089961
$ws.Range("C22").Value = 3
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 7
$ws.Range("G22").Value = 11
$ws.Range("H22").Value = -36.363636363636
$ws.Range("I22").Value = 72
$ws.Range("J22").Value = 50
$ws.Range("K22").Value = 44
$ws.Range("L22").Value = 105.714285714286
$ws.Range("M22").Value = 22.033898305084
$ws.Range("J14").Copy($ws.Range("C23"))
$ws.Range("C23").Value = 1
$ws.Range("G14").Copy($ws.Range("D23"))
$ws.Range("H14").Copy($ws.Range("E23"))
$ws.Range("F23").Value = 5
$ws.Range("H23").Value = 150
$ws.Range("I23").Value = 49
$ws.Range("K23").Value = 28.947368421052
$ws.Range("L23").Value = 6.521739130434
$ws.Range("M23").Value = 28.947368421052
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = -13.043478260869
$ws.Range("F24").Value = 99
$ws.Range("G24").Value = 63
$ws.Range("H24").Value = 57.142857142857
$ws.Range("I24").Value = 924
$ws.Range("J24").Value = 625
$ws.Range("K24").Value = 47.84
$ws.Range("L24").Value = 48.076923076923
$ws.Range("M24").Value = -21.361702127659
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -22.222222222222
$ws.Range("F25").Value = 27
$ws.Range("G25").Value = 28
$ws.Range("H25").Value = -3.571428571428
$ws.Range("I25").Value = 367
$ws.Range("J25").Value = 298
$ws.Range("K25").Value = 23.154362416107
$ws.Range("L25").Value = 37.969924812030
$ws.Range("M25").Value = 38.490566037735
$ws.Range("C26").Value = 2
$ws.Range("J14").Copy($ws.Range("D26"))
$ws.Range("D26").Value = 1
$ws.Range("H15").Copy($ws.Range("E26"))
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 3
$ws.Range("H26").Value = 50
$ws.Range("I26").Value = 31
$ws.Range("J26").Value = 24
$ws.Range("K26").Value = 29.166666666666
$ws.Range("L26").Value = 34.782608695652
$ws.Range("C27").Value = 5
$ws.Range("G14").Copy($ws.Range("D27"))
$ws.Range("H14").Copy($ws.Range("E27"))
$ws.Range("F27").Value = 11
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 83.333333333333
$ws.Range("I27").Value = 103
$ws.Range("K27").Value = -0.961538461538
$ws.Range("L27").Value = 106
$ws.Range("J14").Copy($ws.Range("C28"))
$ws.Range("C28").Value = 1
$ws.Range("J14").Copy($ws.Range("F28"))
$ws.Range("F28").Value = 1
$ws.Range("I28").Value = 7
$ws.Range("K28").Value = 250
$ws.Range("L28").Value = 75
$ws.Range("M28").Value = 600
$ws.Range("N28").Value = -68.181818181818
$ws.Range("J14").Copy($ws.Range("C29"))
$ws.Range("C29").Value = 1
$ws.Range("J14").Copy($ws.Range("F29"))
$ws.Range("F29").Value = 1
$ws.Range("I29").Value = 7
$ws.Range("K29").Value = 250
$ws.Range("L29").Value = 133.333333333333
$ws.Range("M29").Value = 600
$ws.Range("N29").Value = -66.666666666666
$ws.Range("J14").Copy($ws.Range("C30"))
$ws.Range("C30").Value = 1
$ws.Range("I30").Value = 11
$ws.Range("K30").Value = 37.5
$ws.Range("L30").Value = 120
